# Leave Card update - 12/27/2023 4:01 PM
# Adds a 3-day Forced Leave (FL) entry for 12/27-29/2023 between the
# Nov-2023 and Jan-2024 monthly rows, plus two SL(1-0-0) half-day credits
# that were filed against the Oct-2023 and Nov-2023 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Insert a new blank table row at row 406 (pushes rows 406.. down by one) ---
$ws.Rows.Item(406).Insert()
$lo.Resize($ws.Range("A8:K523"))

# Copy the row formatting (borders/number formats/styles) from the row
# directly above (405) into the freshly inserted, still-blank row 406 so it
# matches the rest of the table. PasteSpecial formats-only leaves the new
# row's values/formulas empty.
$ws.Range("A405:K405").Copy()
$ws.Range("A406:K406").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The table Resize() re-derives the calculated column formula for the new
# last row (523) using "[@EARNED]" shorthand, which this engine evaluates
# to a stale #VALUE! cache. Re-assert the original unambiguous formula so
# it recalculates cleanly (matches every other row in the column).
$ws.Range("G523").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]])," + '""' + ",Table1[[#This Row],[EARNED]])"

# --- 2. Row 403 (period ending 9/1/2023->... row for Sep 2023): record a
#        0.75-day SL credit already reflected in EARNED (no particulars text) ---
$ws.Range("C403").Value = 1.25

# --- 3. Row 404: SL(1-0-0) filed, reflected as 1 day Absence Undertime W/Pay,
#        filed/approved 10/13/2023 ---
$ws.Range("B404").Value = "SL(1-0-0)"
$ws.Range("C404").Value = 1.25
$ws.Range("H404").Value = 1
$ws.Range("K397").Copy()
$ws.Range("K404").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K404").Value = 45212

# --- 4. Row 405: SL(1-0-0) filed, reflected as 1 day Absence Undertime W/Pay,
#        filed/approved 11/20/2023 ---
$ws.Range("B405").Value = "SL(1-0-0)"
$ws.Range("C405").Value = 1.25
$ws.Range("H405").Value = 1
$ws.Range("K397").Copy()
$ws.Range("K405").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K405").Value = 45250

# --- 5. New row 406: Forced Leave FL(3-0-0), 3 days Absence Undertime W/O Pay,
#        dated 12/27-29/2023 ---
$ws.Range("B406").Value = "FL(3-0-0)"
$ws.Range("D406").Value = 3
$ws.Range("K399").Copy()
$ws.Range("K406").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K406").Value = "12/27-29/2023"

# --- 6. Recalculate so BALANCE/EARNED totals in the summary rows update ---
$excel.CalculateFullRebuild()

# Reflect the final working selection/scroll position as left by the editor.
$ws.Range("E406").Select()
